$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4:B4").Copy()
$ws.Range("A5:B5").PasteSpecial(-4122)

$ws.Range("A5").Value = -15
$ws.Range("B5").Value = 0

$ws.Range("A6").Select()
